# data/exp1_demographic.xlsx — "Add files via upload"
#
# The sheet holds one demographic row per participant (participant id,
# biological sex, age, years of education). Row 24 (participant "s48")
# was missing its age/education figures — fill them in now that the
# values are available, then leave the window scrolled/selected the way
# it was when the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24: age = 41, years of education = 18
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 18

# Restore the on-screen view state captured in the saved file: scrolled
# so row 10 is the first visible row, with J20 as the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J20").Select()
